$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells whose new values look like plain numbers need to be forced to
# text format first, matching the original inline-string cell type,
# otherwise Excel will coerce them into numeric cells.
$textForceRefs = @("D5", "D6", "D7", "D9", "D10", "D11", "D19", "D20", "D23", "D24", "D25", "D27", "D28", "D29", "D31", "D32", "D34", "D40", "D43", "D46", "D49", "D51")
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).NumberFormat = "@"
}

# Apply the updated crypto price / volume values
$ws.Range("D2").Value = "47.936.25"
$ws.Range("E2").Value = "  -0.43%  "
$ws.Range("D3").Value = "2.481.42"
$ws.Range("E3").Value = "  -1.48%  "
$ws.Range("E4").Value = "  -0.07%  "
$ws.Range("D5").Value = "317.08"
$ws.Range("E5").Value = "  -1.55%  "
$ws.Range("D6").Value = "104.73"
$ws.Range("E6").Value = "  -4.46%  "
$ws.Range("D7").Value = "0.518"
$ws.Range("E7").Value = "  -2.76%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").Value = "0.535"
$ws.Range("E9").Value = "  -3.57%  "
$ws.Range("D10").Value = "38.92"
$ws.Range("E10").Value = "  -4.09%  "
$ws.Range("D11").Value = "20.37"
$ws.Range("E11").Value = "  -0.93%  "
$ws.Range("E12").Value = "  -3.24%  "
$ws.Range("E13").Value = "  +0.31%  "
$ws.Range("E14").Value = "  -3.31%  "
$ws.Range("D15").Value = "2.868.10"
$ws.Range("E15").Value = "  -1.61%  "
$ws.Range("D16").Value = "2.486.32"
$ws.Range("E16").Value = "  -1.40%  "
$ws.Range("E17").Value = "  -3.73%  "
$ws.Range("D18").Value = "47.863.83"
$ws.Range("E18").Value = "  -0.25%  "
$ws.Range("D19").Value = "2.96"
$ws.Range("E19").Value = "  +9.68%  "
$ws.Range("D20").Value = "12.66"
$ws.Range("E20").Value = "  -4.18%  "
$ws.Range("E21").Value = "  -1.93%  "
$ws.Range("D22").Value = "0.0₃0926"
$ws.Range("E22").Value = "  -2.30%  "
$ws.Range("D23").Value = "276.64"
$ws.Range("E23").Value = "  +4.60%  "
$ws.Range("D24").Value = "70.67"
$ws.Range("E24").Value = "  -2.03%  "
$ws.Range("D25").Value = "2.49"
$ws.Range("E25").Value = "  -3.38%  "
$ws.Range("E26").Value = "  +0.23%  "
$ws.Range("D27").Value = "25.59"
$ws.Range("E27").Value = "  -1.98%  "
$ws.Range("D28").Value = "2.19"
$ws.Range("E28").Value = "  -0.52%  "
$ws.Range("D29").Value = "9.58"
$ws.Range("E29").Value = "  -5.52%  "
$ws.Range("E30").Value = "  -5.40%  "
$ws.Range("D31").Value = "34.58"
$ws.Range("E31").Value = "  -5.00%  "
$ws.Range("D32").Value = "49.19"
$ws.Range("E32").Value = "  -1.12%  "
$ws.Range("E33").Value = "  -0.19%  "
$ws.Range("D34").Value = "18.93"
$ws.Range("E34").Value = "  -4.61%  "
$ws.Range("E35").Value = "  -3.05%  "
$ws.Range("E36").Value = "  -2.82%  "
$ws.Range("E37").Value = "  -2.59%  "
$ws.Range("E38").Value = "  -4.37%  "
$ws.Range("E39").Value = "  -4.87%  "
$ws.Range("D40").Value = "122.27"
$ws.Range("E40").Value = "  +1.96%  "
$ws.Range("E41").Value = "  -1.69%  "
$ws.Range("E42").Value = "  +0.50%  "
$ws.Range("D43").Value = "21.98"
$ws.Range("E43").Value = "  -0.64%  "
$ws.Range("E44").Value = "  -0.09%  "
$ws.Range("D45").Value = "1.992.57"
$ws.Range("E45").Value = "  -1.21%  "
$ws.Range("D46").Value = "3.13"
$ws.Range("E46").Value = "  -1.14%  "
$ws.Range("E47").Value = "  -1.66%  "
$ws.Range("E48").Value = "  -3.20%  "
$ws.Range("D49").Value = "8.90"
$ws.Range("E49").Value = "  -2.58%  "
$ws.Range("E50").Value = "  -2.32%  "
$ws.Range("D51").Value = "78.75"
$ws.Range("E51").Value = "  -0.20%  "

# Restore default (unstyled) formatting on the cells we forced to text,
# so the cell style matches the original (no explicit style index).
foreach ($ref in $textForceRefs) {
    $ws.Range($ref).Style = "Normal"
}
